$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$row = 7
foreach ($goal in $goals) {
    $ws.Cells.Item($row, 1).Value = $goal.Id
    $ws.Cells.Item($row, 2).Value = $goal.Name
    $ws.Cells.Item($row, 3).Value = 45908
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value = 0.9900990099009901
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = -0.01
    $row++
}
